$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.546.69"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.812.18"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.96"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.579"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.92"
$ws.Range("E8").Value = "  +6.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.301"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "2.074.04"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.23"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.820.46"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.651"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.532.94"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.46"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.28"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "0.0₃0799"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.51"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "172.67"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.98"
$ws.Range("E26").Value = "  +9.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.84"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.04"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "1.397.39"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("E37").Value = "  -5.72%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.80"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.963"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.84"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.40"
$ws.Range("E44").Value = "  -3.38%  "
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0516"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").Value = "1.973.68"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.32"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("E51").Value = "  +0.18%  "

Write-Host "Cryptos list updated."
